$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 142, shifting existing rows 142:163 down to 143:164.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new record.
$ws.Range("A142").Value = 3
$ws.Range("B142").Value = 'Femacal de La Calera'
$ws.Range("C142").Value = 'Coquimbo'
$ws.Range("D142").Value = 44491
$ws.Range("E142").Value = 5
$ws.Range("F142").Value = 100112001
$ws.Range("G142").Value = 'Berenjena'
$ws.Range("H142").Value = 'Sin especificar'
$ws.Range("I142").Value = 'Primera'
$ws.Range("J142").Value = 105
$ws.Range("K142").Value = 7500
$ws.Range("L142").Value = 8000
$ws.Range("M142").Value = 7762
$ws.Range("N142").Value = '$/caja 60 unidades'
$ws.Range("O142").Value = 'Región de Arica y Parinacota'
$ws.Range("P142").Value = 129
$ws.Range("Q142").Value = 60
$ws.Range("R142").Value = 'Hortaliza'
